# "fabric labeliling module refactoring"
#
# - sheet "report": move selection from C5 -> B9
# - sheet "service_tables": move selection from E75 -> H41 (and scroll to E19,
#   best-effort — the engine does not persist topLeftCell)
# - sheet "service_tables": flip several keys/export_to_excel flags (cols G/H)
# - sheet "service_tables": remove the AutoFilter

$wb = $excel.ActiveWorkbook

# --- sheet "report": selection C5 -> B9 -------------------------------------
$wsReport = $wb.Worksheets.Item("report")
$wsReport.Activate()
$wsReport.Range("B9").Select()

# --- sheet "service_tables": data + view changes -----------------------------
$wsSvc = $wb.Worksheets.Item("service_tables")
$wsSvc.Activate()

# keys / export_to_excel flag flips
$wsSvc.Range("H36").Value = 1

$wsSvc.Range("G62").Value = 0

$wsSvc.Range("H70").Value = 0
$wsSvc.Range("G71").Value = 0
$wsSvc.Range("G72").Value = 0
$wsSvc.Range("G73").Value = 0
$wsSvc.Range("G74").Value = 0
$wsSvc.Range("G75").Value = 0
$wsSvc.Range("G76").Value = 0
$wsSvc.Range("G77").Value = 0
$wsSvc.Range("G78").Value = 0
$wsSvc.Range("G81").Value = 0

# remove the AutoFilter that used to cover A1:J81
$wsSvc.AutoFilterMode = $false

# scroll position (topLeftCell C55 -> E19) then land the selection on H41
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 5
$wsSvc.Range("H41").Select()
